$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain decimal number and must be
# forced to Text format so Excel does not auto-convert it to a Number,
# matching the source data which stores these as inline strings.
$textForceCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "26.042.26"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.647.57"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "206.60"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "0.5193"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.2579"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "0.06254"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "20.72"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "0.07555"
$ws.Range("D12").Value = "1.647.92"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "4.378"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").Value = "66.13"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "0.0₅7930"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "26.080.45"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "4.666"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "187.62"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "10.01"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").Value = "6.130"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "148.09"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "0.1212"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "7.352"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").Value = "15.64"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "1.388"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").Value = "0.06015"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("D30").Value = "1.242"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").Value = "3.441"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "3.392"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "1.624"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "0.9764"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "2.381"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "2.736"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "0.5886"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").Value = "0.01592"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "5.963"
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.079.24"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.8452"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "100.25"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").Value = "1.802.14"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "0.0₈107"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "1.009"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "54.71"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").Value = "7.998"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "0.05217"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "0.4238"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "5.863"
$ws.Range("E51").Value = "  -0.39%  "
